$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match style of existing header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells F2:F7 - plain text timestamps, same formatting as data cells in column E
$timestamps = @(
    "2021-10-05 13:39:05.069746",
    "2021-10-05 13:39:05.069760",
    "2021-10-05 13:39:05.069764",
    "2021-10-05 13:39:05.069769",
    "2021-10-05 13:39:05.069773",
    "2021-10-05 13:39:05.069776"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
